# tdf#141309 tdf#142215 style fixture update:
# Add a new data row (time value 0.5 / "c") to the autofilter range and
# extend the autofilter (and its backing _FilterDatabase defined name) to
# include it, with an extra discrete filter value for 0.500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the current autofilter so the range/criteria can be redefined cleanly.
$ws.AutoFilterMode = $false

# New row 8: A8 is a time-of-day value (12:00 -> 0.5 as a day fraction),
# B8 reuses the existing shared string "c".
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Value = "c"

# Re-apply the autofilter over the expanded A1:B8 range, filtering column A
# (field 1) on the discrete values 0.046, 0.500 and 0.516 (xlFilterValues = 7).
$ws.Range("A1:B8").AutoFilter(1, @("0.046", "0.500", "0.516"), 7)

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the new range.
foreach ($n in $ws.Names) {
    if ($n.Name -eq "Munka1!_FilterDatabase") {
        $n.RefersTo = "=Munka1!`$A`$1:`$B`$8"
    }
}

# Match the author's cursor position after the edit.
$ws.Range("C7").Select()
